$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Drop the oldest quarter column (D) so every later quarter shifts left ---
$ws.Columns("D").Delete()

# --- Append a new quarter column at M (copy L's formatting, including column width) ---
$ws.Range("L1:L28").Copy($ws.Range("M1"))
$ws.Columns("M").ColumnWidth = 30.17

# --- Fix header row 8 (quarter label) and row 9 (publish date) ---
$ws.Range("I9").Value = "1402-02-28 (8)"
$ws.Range("M8").Value = "فصل چهارم منتهی به 1401/12"
$ws.Range("M9").Value = "1402-02-28"

# --- Fix the two cells whose data didn't simply shift (I25, I26) ---
$ws.Range("I25").Value = 0
$ws.Range("I26").Value = 65351

# --- New quarter's financial figures (column M) ---
$ws.Range("M11").Value = 31939
$ws.Range("M12").Value = -24347
$ws.Range("M13").Value = 7592
$ws.Range("M14").Value = -2378
$ws.Range("M16").Value = 1910
$ws.Range("M17").Value = 7124
$ws.Range("M18").Value = -601
$ws.Range("M19").Value = 13
$ws.Range("M20").Value = 6536
$ws.Range("M21").Value = 1923
$ws.Range("M22").Value = 8459
$ws.Range("M24").Value = 8459
$ws.Range("M25").Value = 0
$ws.Range("M26").Value = 39601
